# primer commit sistema de legalizacion
#
# Renumbers the example invoice/delivery/order numbers from the
# 2025-0106..0108 batch to the 2025-0109..0111 batch, mirrors the invoice
# number into the delivery/internal-order/purchase-order columns on the
# Headers sheet, and leaves the UI focused on the Headers sheet (matching
# the author's last on-screen selection) instead of Details.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Headers")
$ws2 = $wb.Worksheets.Item("Details")

# --- Renumber the shared invoice/delivery/order numbers everywhere they
#     occur (Headers!A2:A4 and Details!A2:A10) -----------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("2025-0106", "2025-0109")
    $ws.Cells.Replace("2025-0107", "2025-0110")
    $ws.Cells.Replace("2025-0108", "2025-0111")
}

# --- Mirror the (now renumbered) invoice number into the delivery number,
#     internal order number and purchase order number columns on Headers --
$ws1.Range("B2").Value = "2025-0109"
$ws1.Range("C2").Value = "2025-0109"
$ws1.Range("D2").Value = "2025-0109"

$ws1.Range("B3").Value = "2025-0110"
$ws1.Range("C3").Value = "2025-0110"
$ws1.Range("D3").Value = "2025-0110"

$ws1.Range("B4").Value = "2025-0111"
$ws1.Range("C4").Value = "2025-0111"
$ws1.Range("D4").Value = "2025-0111"

# --- Selection / active-sheet bookkeeping ---------------------------------
# Leave Details with a selection on A8:A10 ...
$ws2.Range("A8:A10").Select()
# ... then switch to, and leave the workbook focused on, Headers with D10
# selected (matches the saved view state captured in the workbook).
$ws1.Activate()
$ws1.Range("D10").Select()
